$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing row values (now representing ECs -> Efna4/Epha5 -> sCs)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna4"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6731353333333333
$ws.Range("H2").Value = 2.019406
$ws.Range("I2").Value = 0.3272865828458516
$ws.Range("J2").Value = 0.3272865828458516
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1979113333333334
$ws.Range("N2").Value = 0.5937340000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1332211113337778
$ws.Range("R2").Value = 1.198990002004
$ws.Range("S2").Value = 0.3272865828458516
$ws.Range("T2").Value = 0.3272865828458516

# Row 3: update existing row values (FAPs -> Efna4/Epha5 -> sCs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna4"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9964423333333334
$ws.Range("H3").Value = 2.989327
$ws.Range("I3").Value = 0.4844823769162027
$ws.Range("J3").Value = 0.4844823769162026
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1979113333333334
$ws.Range("N3").Value = 0.5937340000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1972072307797778
$ws.Range("R3").Value = 1.774865077018001
$ws.Range("S3").Value = 0.4844823769162027
$ws.Range("T3").Value = 0.4844823769162026

# Row 4: new row (sCs -> Efna4/Epha5 -> sCs)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Efna4"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3871376666666667
$ws.Range("H4").Value = 1.161413
$ws.Range("I4").Value = 0.1882310402379457
$ws.Range("J4").Value = 0.1882310402379457
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1979113333333334
$ws.Range("N4").Value = 0.5937340000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.07661893179355556
$ws.Range("R4").Value = 0.6895703861420002
$ws.Range("S4").Value = 0.1882310402379457
$ws.Range("T4").Value = 0.1882310402379457
